# Update Betfair Back/Lay odds for 2026-01-06 per latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Australian A-League Men: Melbourne City vs Brisbane Roar
$ws.Range("K2").Value = 3.8
$ws.Range("T2").Value = 2
$ws.Range("AG2").Value = 10.5
$ws.Range("AI2").Value = 90
$ws.Range("AO2").Value = 100

# Row 3 - Italian Serie A: Pisa vs Como
$ws.Range("F3").Value = 5.6
$ws.Range("G3").Value = 5.8
$ws.Range("H3").Value = 1.8
$ws.Range("J3").Value = 3.7
$ws.Range("N3").Value = 3.3
$ws.Range("O3").Value = 1.41
$ws.Range("Q3").Value = 2.22
$ws.Range("U3").Value = 1.89
$ws.Range("W3").Value = 1.21
$ws.Range("Y3").Value = 7.4
$ws.Range("AB3").Value = 16.5
$ws.Range("AL3").Value = 95
$ws.Range("AO3").Value = 14

# Row 4
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = 6.2
$ws.Range("H4").Value = 1.8
$ws.Range("I4").Value = 1.81
$ws.Range("J4").Value = 3.5
$ws.Range("K4").Value = 3.55
$ws.Range("N4").Value = 3
$ws.Range("O4").Value = 1.47
$ws.Range("T4").Value = 2.26
$ws.Range("U4").Value = 1.76
$ws.Range("V4").Value = 2.22
$ws.Range("W4").Value = 1.19
$ws.Range("Z4").Value = 9
$ws.Range("AA4").Value = 18
$ws.Range("AE4").Value = 23
$ws.Range("AF4").Value = 42
$ws.Range("AG4").Value = 24
$ws.Range("AH4").Value = 26
$ws.Range("AK4").Value = 110
$ws.Range("AL4").Value = 130
$ws.Range("AM4").Value = 210

# Row 5
$ws.Range("F5").Value = 5.2
$ws.Range("H5").Value = 1.79
$ws.Range("I5").Value = 1.8
$ws.Range("J5").Value = 3.95
$ws.Range("K5").Value = 4
$ws.Range("O5").Value = 1.36
$ws.Range("Q5").Value = 2.04
$ws.Range("T5").Value = 1.97
$ws.Range("U5").Value = 1.98
$ws.Range("V5").Value = 2.24
$ws.Range("X5").Value = 13.5
$ws.Range("Y5").Value = 8.199999999999999
$ws.Range("AC5").Value = 8.4
$ws.Range("AL5").Value = 80

# Row 6
$ws.Range("F6").Value = 2.82
$ws.Range("G6").Value = 2.96
$ws.Range("H6").Value = 2.64
$ws.Range("I6").Value = 2.74
$ws.Range("N6").Value = 3.55
$ws.Range("P6").Value = 1.86
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 1.31
$ws.Range("S6").Value = 3.55
$ws.Range("T6").Value = 1.75
$ws.Range("U6").Value = 2.12
$ws.Range("V6").Value = 1.59
$ws.Range("W6").Value = 1.5
$ws.Range("X6").Value = 14
$ws.Range("Z6").Value = 17.5
$ws.Range("AA6").Value = 40
$ws.Range("AB6").Value = 11.5
$ws.Range("AH6").Value = 21
$ws.Range("AI6").Value = 46
$ws.Range("AJ6").Value = 46
$ws.Range("AK6").Value = 34
$ws.Range("AL6").Value = 48
$ws.Range("AM6").Value = 100
$ws.Range("AN6").Value = 32

# Row 7
$ws.Range("L7").Value = 1.42
$ws.Range("R7").Value = 1.37
$ws.Range("U7").Value = 2.2
$ws.Range("X7").Value = 13.5
$ws.Range("AB7").Value = 12.5
$ws.Range("AF7").Value = 21
$ws.Range("AG7").Value = 13.5
$ws.Range("AK7").Value = 36

# Row 8
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = 6.6
$ws.Range("J8").Value = 4.1
$ws.Range("T8").Value = 1.82
$ws.Range("W8").Value = 2.5
$ws.Range("AE8").Value = 90
$ws.Range("AI8").Value = 80
$ws.Range("AJ8").Value = 15.5
